$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3838.793663506989
$ws.Range("G2").Value = 3844.457907582817
$ws.Range("F3").Value = 3240.698776960683
$ws.Range("G3").Value = 3244.860551470487
$ws.Range("F4").Value = 4298.890449645307
$ws.Range("G4").Value = 4302.13548085098
$ws.Range("F5").Value = 4651.307404371489
$ws.Range("G5").Value = 4656.523875577874
$ws.Range("F6").Value = 5081.186898040413
$ws.Range("G6").Value = 5087.073750415043
$ws.Range("F7").Value = 3645.395091249914
$ws.Range("G7").Value = 3649.347854986108
$ws.Range("F8").Value = 2873.518666666828
$ws.Range("G8").Value = 2884.698833333496
$ws.Range("F9").Value = 2400.104805265191
$ws.Range("G9").Value = 2406.905540028984
$ws.Range("F10").Value = 2367.561208875421
$ws.Range("G10").Value = 2373.877677888428
$ws.Range("F11").Value = 2999.752300613666
$ws.Range("G11").Value = 3007.881518405077
$ws.Range("F12").Value = 1354.432470024058
$ws.Range("G12").Value = 1358.055497601995
$ws.Range("F13").Value = 2076.274270918428
$ws.Range("G13").Value = 2079.752996968478
$ws.Range("F14").Value = 1015.763804401034
$ws.Range("G14").Value = 1019.669068459715
$ws.Range("F15").Value = 957.6159940653357
$ws.Range("G15").Value = 961.9051793903969
$ws.Range("F16").Value = 265.0762463343253
$ws.Range("G16").Value = 265.3610703812464
$ws.Range("F17").Value = 2840.563752784126
$ws.Range("G17").Value = 2850.56827672622
$ws.Range("F18").Value = 1619.359411764797
$ws.Range("G18").Value = 1627.89117647068
$ws.Range("F19").Value = 1198.500707018038
$ws.Range("G19").Value = 1202.127634713521
